# Update the "NumPoliza" data value in the Sura "Endoso - Adicional" DataSource sheet
# and leave the selection where the user last clicked (J12), matching the
# Ranorex DataSource reconfiguration described in the commit message.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# E2 holds the policy number under header "NumPoliza" (E1). Keep it as text
# (leading zeros matter) by prefixing with an apostrophe, which preserves the
# existing quoted-text cell style instead of allocating a new one.
$ws.Range("E2").Value = "'04104013002"

# Move/leave the active selection on J12, matching the saved worksheet view.
$ws.Range("J12").Select() | Out-Null
